# Auto-generated: apply crypto price/volume updates from GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "43.421.31"
Set-TextCell "E2" "  -1.52%  "
Set-TextCell "D3" "2.257.02"
Set-TextCell "E3" "  -0.05%  "
Set-TextCell "E4" "  -0.14%  "
Set-TextCell "D5" "231.54"
Set-TextCell "E5" "  +0.78%  "
Set-TextCell "E6" "  +1.37%  "
Set-TextCell "D7" "64.43"
Set-TextCell "E7" "  +1.78%  "
Set-TextCell "E8" "  -0.15%  "
Set-TextCell "D9" "0.439"
Set-TextCell "E9" "  -1.32%  "
Set-TextCell "D10" "0.0955"
Set-TextCell "E10" "  -7.28%  "
Set-TextCell "E11" "  -0.01%  "
Set-TextCell "D12" "26.35"
Set-TextCell "E12" "  +0.51%  "
Set-TextCell "E13" "  -1.98%  "
Set-TextCell "D14" "2.592.39"
Set-TextCell "E14" "  -0.18%  "
Set-TextCell "D15" "14.97"
Set-TextCell "E15" "  -4.21%  "
Set-TextCell "E16" "  -1.38%  "
Set-TextCell "E17" "  -1.04%  "
Set-TextCell "D18" "2.256.96"
Set-TextCell "E18" "  -0.34%  "
Set-TextCell "D19" "43.404.18"
Set-TextCell "E19" "  -1.38%  "
Set-TextCell "D20" "0.0₃0967"
Set-TextCell "E20" "  -4.01%  "
Set-TextCell "D21" "72.93"
Set-TextCell "E21" "  -0.37%  "
Set-TextCell "E22" "  +1.24%  "
Set-TextCell "D23" "247.26"
Set-TextCell "E23" "  -1.46%  "
Set-TextCell "E24" "  +19.24%  "
Set-TextCell "E25" "  -0.02%  "
Set-TextCell "D26" "2.42"
Set-TextCell "E26" "  +0.36%  "
Set-TextCell "E27" "  -2.18%  "
Set-TextCell "D28" "9.71"
Set-TextCell "E28" "  -2.56%  "
Set-TextCell "D29" "173.88"
Set-TextCell "E29" "  +1.17%  "
Set-TextCell "D30" "21.65"
Set-TextCell "E30" "  +4.60%  "
Set-TextCell "E31" "  +3.69%  "
Set-TextCell "E32" "  -4.37%  "
Set-TextCell "E33" "  +0.94%  "
Set-TextCell "E34" "  +4.74%  "
Set-TextCell "D35" "0.0678"
Set-TextCell "E35" "  -0.21%  "
Set-TextCell "E36" "  +1.12%  "
Set-TextCell "D37" "3.61"
Set-TextCell "E37" "  -4.75%  "
Set-TextCell "E38" "  -3.17%  "
Set-TextCell "D39" "2.26"
Set-TextCell "E39" "  -1.52%  "
Set-TextCell "E40" "  -2.64%  "
Set-TextCell "D41" "1.00"
Set-TextCell "E41" "  -0.20%  "
Set-TextCell "D42" "8.79"
Set-TextCell "E42" "  +5.99%  "
Set-TextCell "D43" "4.51"
Set-TextCell "E43" "  +4.17%  "
Set-TextCell "D44" "17.15"
Set-TextCell "E44" "  -1.44%  "
Set-TextCell "D45" "96.59"
Set-TextCell "E45" "  -0.71%  "
Set-TextCell "B46" "Celestia"
Set-TextCell "C46" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D46" "10.15"
Set-TextCell "E46" "  +4.76%  "
Set-TextCell "E47" "  -0.53%  "
Set-TextCell "B48" "Cronos"
Set-TextCell "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D48" "0.0938"
Set-TextCell "E48" "  -2.47%  "
Set-TextCell "E49" "  -1.30%  "
Set-TextCell "D50" "1.425.84"
Set-TextCell "E50" "  -0.88%  "
Set-TextCell "E51" "  -0.73%  "
